$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.954.60"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "2.681.83"
$ws.Range("E3").Value = "  -3.03%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'550.34"
$ws.Range("E5").Value = "  -4.81%  "
$ws.Range("D6").Value = "'157.27"
$ws.Range("E6").Value = "  -2.33%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").Value = "'0.106"
$ws.Range("E9").Value = "  -4.78%  "
$ws.Range("E10").Value = "  -2.92%  "
$ws.Range("E11").Value = "  -5.11%  "
$ws.Range("E12").Value = "  -12.85%  "
$ws.Range("D13").Value = "3.157.89"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").Value = "'25.99"
$ws.Range("E14").Value = "  -5.18%  "
$ws.Range("D15").Value = "62.828.23"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "'0.0000147"
$ws.Range("E16").Value = "  -3.77%  "
$ws.Range("D17").Value = "2.683.56"
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("E19").Value = "  -6.03%  "
$ws.Range("D20").Value = "'342.84"
$ws.Range("E20").Value = "  -4.55%  "
$ws.Range("D21").Value = "'6.30"
$ws.Range("E21").Value = "  -5.74%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'0.504"
$ws.Range("E23").Value = "  -5.21%  "
$ws.Range("D24").Value = "'63.47"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").Value = "'0.168"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "'8.13"
$ws.Range("E27").Value = "  -5.97%  "
$ws.Range("D28").Value = "0.0₃0856"
$ws.Range("E28").Value = "  -8.09%  "
$ws.Range("E29").Value = "  -2.79%  "
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("D31").Value = "'7.02"
$ws.Range("E31").Value = "  -5.06%  "
$ws.Range("D32").Value = "'165.22"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'4.82"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("D35").Value = "'19.54"
$ws.Range("D36").Value = "'1.43"
$ws.Range("E36").Value = "  -6.25%  "
$ws.Range("E37").Value = "  -4.27%  "
$ws.Range("D38").Value = "'338.81"
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("D39").Value = "'6.16"
$ws.Range("E39").Value = "  -4.52%  "
$ws.Range("D40").Value = "'0.931"
$ws.Range("E40").Value = "  -7.79%  "
$ws.Range("D41").Value = "'38.22"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "'3.94"
$ws.Range("E42").Value = "  -6.34%  "
$ws.Range("D43").Value = "'20.36"
$ws.Range("E43").Value = "  -5.84%  "
$ws.Range("D44").Value = "'20.77"
$ws.Range("E44").Value = "  -8.05%  "
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").Value = "'0.0560"
$ws.Range("E46").Value = "  -6.01%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "'0.0973"
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("D50").Value = "'129.50"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("D51").Value = "2.088.77"
$ws.Range("E51").Value = "  -3.07%  "
